# Banco_Dashboard.xlsx - mark the two "tanques de lama" BV activities
# (rows 49-50, sheet "dados_corrigidos") as fully completed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing date-formatted style (numFmtId 14, short date) from F2
# onto F49:G50 so the new date cells reuse the workbook's existing style
# instead of Excel fabricating a brand-new number format.
$ws.Range("F2").Copy()
$ws.Range("F49:G50").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Row 49 - "Abrir BV tanques de lama e slaker A e B"
$ws.Range("D49").Value = 1
$ws.Range("F49").Value = "8/10/2025"
$ws.Range("G49").Value = "8/10/2025"
$ws.Range("H49").Value = 100

# Row 50 - "Abrir BV tanques de lama Slurry A e B"
$ws.Range("D50").Value = 1
$ws.Range("F50").Value = "8/10/2025"
$ws.Range("G50").Value = "8/10/2025"
$ws.Range("H50").Value = 100

# Leave the active selection on D2, matching the saved workbook state.
$ws.Range("D2").Select()
